$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: account number cleared (cell kept, just emptied), name/pay updated
$ws.Range("A2").ClearContents()
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "fgfg fgfgfg"
$ws.Range("C2").Value = 0.86
$ws.Range("D2").Value = "November"

# Row 3 is removed entirely (shifts used range/dimension back to A1:D2)
$ws.Rows.Item(3).Delete()

# Column width tweaks (A and C got narrower) -- ColumnWidth is quantized to
# whole pixels by Excel, so these character-width inputs land on the closest
# achievable pixel boundary to the target stored widths (14.0899.. / 7.4899..)
$ws.Columns.Item(1).ColumnWidth = 13.428571428571429
$ws.Columns.Item(3).ColumnWidth = 6.714285714285714
